# Insert a new data row above the current row 192 (shifting existing rows
# 192-282 down to 193-283) and populate it with a new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(192).Insert()

$ws.Range("A192").Value = 6
$ws.Range("B192").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C192").Value = "Metropolitana"
$ws.Range("D192").Value = 44489
$ws.Range("E192").Value = 13
$ws.Range("F192").Value = 100112043
$ws.Range("G192").Value = "Pepino ensalada"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 610
$ws.Range("K192").Value = 6000
$ws.Range("L192").Value = 7000
$ws.Range("M192").Value = 6426
$ws.Range("N192").Value = "`$/caja 60 unidades"
$ws.Range("O192").Value = "Región de Arica y Parinacota"
$ws.Range("P192").Value = 107
$ws.Range("Q192").Value = 60
$ws.Range("R192").Value = "Hortaliza"
